# xl2yml.py: streamlined generation of C, K and F matrices
#
# The old "elements" sheet is dropped; "elements_new" becomes the new
# "elements" sheet. Its Capacity column (C) is reformatted to show two
# decimal places, and it becomes the active/selected sheet (with D31
# selected), while "flows" keeps its own selection but is no longer the
# active tab.

$wb = $excel.ActiveWorkbook

# 1. Drop the old "elements" sheet entirely.
$oldElements = $wb.Worksheets.Item("elements")
$oldElements.Delete() | Out-Null

# 2. Promote "elements_new" to be the "elements" sheet.
$elements = $wb.Worksheets.Item("elements_new")
$elements.Name = "elements"

# 3. Reformat the Capacity column (C) to 2 decimal places for the whole
#    used range (header + 11 data rows).
$elements.Range("C1:C12").NumberFormat = "0.00"

# Give it an explicit page setup (A4 portrait) like the other sheet.
$ps = $elements.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# 4. Make "elements" the active sheet/tab, with D31 selected.
$elements.Activate() | Out-Null
$elements.Range("D31").Select() | Out-Null

# "flows" retains its own selection (D2) from before.
$flows = $wb.Worksheets.Item("flows")
$flows.Range("D2").Select() | Out-Null
$elements.Activate() | Out-Null
